$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'high waist capris',
    'high waist cigarette pants',
    'high waist cigarette pants for women',
    'high waist compression capri',
    'high waist compression capri leggings',
    'high waist compression capris women',
    'high waist compression leggings',
    'high waist compression leggings for women',
    'high waist compression leggings women',
    'high waist compression shorts women',
    'high waist compression tights women',
    'high waist exercise pants for women',
    'high waist hosen',
    'high waist insulated leggings',
    'high waist leggingd',
    'high waist leggings',
    'high waist leggings black',
    'high waist leggings compression',
    'high waist leggings for women',
    'high waist leggings thick',
    'high waist pants',
    'high waist pants for women',
    'high waist running leggings',
    'high waist running pants',
    'high waist running pants women',
    'high waist ski pants women',
    'high waist tights',
    'high waist tights for women',
    'high waist training leggings',
    'high waist training leggings women',
    'high waist warm leggings',
    'high waist workout',
    'high waist workout capris',
    'high waist yoga',
    'high waist yoga shorts',
    'high waistband leggings',
    'high waisted adidas leggings women',
    'high waisted bermuda',
    'high waisted black pants women',
    'high waisted compression capris',
    'high waisted compression leggings for women',
    'high waisted leggings colorful',
    'high waisted pant',
    'high waisted running shorts',
    'high waisted tights plus size',
    'high waisted tights women',
    'high waisted workout leggings for women plus size',
    'high weist lycra leggins',
    'high yoga pants for women',
    'high-performance neck fan',
    'hight waist pants women',
    'hight waist tights',
    'hike flash light',
    'hike pants women',
    'hike sticks for women',
    'hiking capri men',
    'hiking capri women',
    'hiking capris for women',
    'hiking care package',
    'hiking cold weather pants women',
    'hiking hip pack',
    'hiking knee brace',
    'hiking knee support',
    'hiking leggings',
    'hiking leggings for women',
    'hiking leggings women',
    'hiking panties for women',
    'hiking tights',
    'hinge bike',
    'hinge clothing for women',
    'hinge tank top',
    'hinge womens clothing',
    'hip and joint for women',
    'hip and joint pro',
    'hip and joint supplements for women',
    'hip and joint xl',
    'hip art',
    'hip boys jacket',
    'hip brace for women',
    'hip brace for women support',
    'hip brace women',
    'hip braces for women',
    'hip clothing for women',
    'hip compression',
    'hip compression brace',
    'hip compression shorts women',
    'hip compression women',
    'hip exercise',
    'hip flexor',
    'hip flexor pain',
    'hip flexor stretch',
    'hip flexor stretches',
    'hip flexor support',
    'hip hip',
    'hip hugger shorts for women',
    'hip joint',
    'hip joint model',
    'hip joint pain brace',
    'hip joint support',
    'hip joint support brace'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
